# Updated DPM integration testfixture with hierarchy node labels
#
# - regenerates the random GUID identifiers in the ID columns across the
#   CodeSchemes / Codes / Extensions / Members_EDA-H1 / Members_EDA-H2 sheets
# - renames "EDA member N" -> "EDA hierarchy node N" labels on the two
#   hierarchy-members sheets
# - prefixes a couple of RELATION codes with "code:"
# - column widths are bumped to reflect the new (generally longer) label text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# CodeSchemes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CodeSchemes")
$ws.Range("A2").Value = "f01bc58e-d294-4d60-a62c-ef0670327a6c"

# ---------------------------------------------------------------------------
# Codes
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Codes")
$ws.Range("A2").Value  = "ef09b5fe-bacb-47c0-95bd-430a8e099807"
$ws.Range("A3").Value  = "a91bbd3c-4eed-449f-a490-6ecd56e414e0"
$ws.Range("A4").Value  = "37d60266-6024-40fa-b243-f7e421575ff4"
$ws.Range("A5").Value  = "ccb929a9-f1d8-43a9-a102-45e5e1fbdb10"
$ws.Range("A6").Value  = "4ffad908-d016-4218-bf39-17889e72b2db"
$ws.Range("A7").Value  = "39145216-77c6-4793-97dc-5a85c96ae4b6"
$ws.Range("A8").Value  = "226a3cd5-2769-4893-a741-95f288aac396"
$ws.Range("A9").Value  = "b2e795f8-a69b-4105-b453-83fd2ebf47ab"
$ws.Range("A10").Value = "90a867d5-eb5e-4b2e-9469-f47005671446"
$ws.Range("A11").Value = "836d14a1-e2af-47b8-8bc3-c6bac35ab91d"
$ws.Range("A12").Value = "41faeb26-20cd-41aa-8cc5-269a7e537156"
$ws.Range("A13").Value = "b4a8397f-3774-4b50-8ae9-fe8482ae2a63"
$ws.Range("A14").Value = "247bd124-b7f7-4676-ae07-b3d3220fb3e5"
$ws.Columns.Item(1).ColumnWidth = 35.57142857142857

# ---------------------------------------------------------------------------
# Extensions
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extensions")
$ws.Range("A2").Value = "149c6a99-b687-492b-a35b-2548ba6c8e65"
$ws.Range("A3").Value = "15db7eaf-9d12-4ab6-bb12-11c87f30efd9"
$ws.Range("A4").Value = "76b41e36-ba71-43c2-84c6-898dacdd843a"
$ws.Columns.Item(1).ColumnWidth = 32.285714285714285

# ---------------------------------------------------------------------------
# Members_EDA-H1
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Members_EDA-H1")
$ws.Range("A2").Value = "6adb97f2-311a-4a78-b2c8-2f2d8bd231bf"
$ws.Range("D2").Value = "EDA hierarchy node 7"

$ws.Range("A3").Value = "3e8b4059-fdb1-4363-a6e6-5736ae2d53a8"
$ws.Range("D3").Value = "EDA hierarchy node 4"

$ws.Range("A4").Value = "9543654f-c9d0-4a93-8363-8f7a9bfe5307"
$ws.Range("C4").Value = "code:EDA-x9"
$ws.Range("D4").Value = "EDA hierarchy node 5"

$ws.Range("A5").Value = "fc457230-0f51-4c8a-b5bc-3b1fdc565d33"
$ws.Range("D5").Value = "EDA hierarchy node 6"

$ws.Range("A6").Value = "a6b92169-eada-474b-ad2d-482b1c205e9c"
$ws.Range("D6").Value = "EDA hierarchy node 2"

$ws.Range("A7").Value = "378d13f9-83c4-4e8d-aa58-09a772e3c8f4"
$ws.Range("C7").Value = "code:EDA-x2"
$ws.Range("D7").Value = "EDA hierarchy node 1"

$ws.Range("A8").Value = "9d318799-0ebf-4278-a415-d1fb7c184569"
$ws.Range("C8").Value = "code:EDA-x2"
$ws.Range("D8").Value = "EDA hierarchy node 3"

$ws.Range("A9").Value = "5d73c619-6ad3-444b-9e76-31cda0fe32dd"
$ws.Range("D9").Value = "EDA hierarchy node 8"

$ws.Columns.Item(1).ColumnWidth = 33.42857142857143
$ws.Columns.Item(4).ColumnWidth = 19.142857142857142

# ---------------------------------------------------------------------------
# Members_EDA-H2
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Members_EDA-H2")
$ws.Range("A2").Value = "dc1465f6-9a1e-4677-a0a7-a2e1d2b151ee"
$ws.Range("F2").Value = "EDA hierarchy node (=, 1)"

$ws.Range("A3").Value = "7feebc59-1457-43c3-9808-fef21f691ae9"
$ws.Range("F3").Value = "EDA hierarchy node (>, 2)"

$ws.Range("A4").Value = "31fc107f-0faa-4bd4-a697-1150a8b580d9"
$ws.Range("F4").Value = "EDA hierarchy node (<, 0)"

$ws.Range("A5").Value = "640aa29f-a992-454e-8b5a-c02b1c88b578"
$ws.Range("F5").Value = "EDA hierarchy node (>=, -1)"

$ws.Range("A6").Value = "48637fb5-2055-4d1d-9963-d49f84c5f897"
$ws.Range("F6").Value = "EDA hierarchy node (<=, -2)"

$ws.Columns.Item(6).ColumnWidth = 23.428571428571427
